# Generate Report for Handoff
# Updates the localization-status report: status moves from "In Translation"
# to "Ready for handoff", and the handoff timestamps are refreshed.

$wb = $excel.ActiveWorkbook
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

# --- Status: "In Translation" -> "Ready for handoff" ---
# Touch every cell that shares this string so the shared-string table ends
# up with a single updated entry instead of a stale duplicate.
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZh.Range("C2").Value = "Ready for handoff"
$wsDe.Range("C2").Value = "Ready for handoff"

# --- Timestamps ---
# Overview!G2 and de-de!H2 shared "2016-08-29 14:43:49" -> bump to 14:44:37
$wsOverview.Range("G2").Value = "2016-08-29 14:44:37"
$wsDe.Range("H2").Value = "2016-08-29 14:44:37"

# zh-cn!H2 held "2016-08-29 14:43:43" -> bump to 14:44:33
$wsZh.Range("H2").Value = "2016-08-29 14:44:33"

# --- Column widths ---
# Widen the Status-related columns (Overview E:F, and column C on the
# per-locale sheets) to fit the new, longer "Ready for handoff" text.
$wsOverview.Columns.Item(5).ColumnWidth = 16.333333333333332
$wsOverview.Columns.Item(6).ColumnWidth = 16.333333333333332
$wsZh.Columns.Item(3).ColumnWidth = 16.333333333333332
$wsDe.Columns.Item(3).ColumnWidth = 16.333333333333332
